$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 445, pushing existing rows 445:476 down to 446:477
$ws.Rows("445:445").Insert()

# Populate the newly inserted row 445 with the new record's data.
# Columns A, B, C, E, F, G, H, R keep the same constant values used throughout
# this sheet subset (Mercado/Región/Codreg/Categoría/Variedad/Clasificación).
$ws.Cells.Item(445, 1).Value = 11
$ws.Cells.Item(445, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(445, 3).Value = "Bíobío"
$ws.Cells.Item(445, 4).Value = 44714
$ws.Cells.Item(445, 5).Value = 8
$ws.Cells.Item(445, 6).Value = 100112004
$ws.Cells.Item(445, 7).Value = "Cebolla"
$ws.Cells.Item(445, 8).Value = "Sin especificar"
$ws.Cells.Item(445, 9).Value = "1a (cosecha)"
$ws.Cells.Item(445, 10).Value = 220
$ws.Cells.Item(445, 11).Value = 7500
$ws.Cells.Item(445, 12).Value = 8000
$ws.Cells.Item(445, 13).Value = 7773
$ws.Cells.Item(445, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(445, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(445, 16).Value = 432
$ws.Cells.Item(445, 17).Value = 18
$ws.Cells.Item(445, 18).Value = "Hortaliza"
